$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The purchase-price columns (I = "Purchase Price", J = "Total Purchase
# Price") hold amounts as text. Force text format first so the new,
# still numeric-looking values ("0.56" etc.) aren't auto-converted to
# numbers when assigned, matching how the original values were stored.
$priceAddrs = @("I2","J2","I3","J3","I4","J4","I5","J5")
foreach ($addr in $priceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 (entry 1, AWS Account): purchase price 0.60 -> 0.56
$ws.Range("I2").Value = "0.56"
$ws.Range("J2").Value = "0.56"

# Row 3 (entry 2, AWS Account): purchase price -0.60 -> -0.56
$ws.Range("I3").Value = "-0.56"
$ws.Range("J3").Value = "-0.56"

# Row 4 (entry 3, GCP Account): purchase price 0.70 -> 0.65
$ws.Range("I4").Value = "0.65"
$ws.Range("J4").Value = "0.65"

# Row 5 (entry 4): purchase price 0.40 -> -0.65, and vendor changes from
# "Azure Account" to "GCP Account" (currency conversion test now only
# exercises AWS/GCP accounts).
$ws.Range("I5").Value = "-0.65"
$ws.Range("J5").Value = "-0.65"

# Restore the default (General) style on the price cells now that the
# text values are safely stored, so no visible formatting change remains.
foreach ($addr in $priceAddrs) {
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("L5").Value = "GCP Account"
$ws.Range("M5").Value = "6b65a6a4-8b81-48f6-b38a-088ca65ed389"
